# Update column G ("K") values on the active sheet for rows 2-24.
# These represent regenerated "K" values (replacing the old "Strike#" derived
# values) after recalculating std/mean and writing s_vals upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 2
    8  = 3
    9  = 1
    10 = 2
    11 = 2
    12 = 3
    13 = 1
    14 = 2
    15 = 2
    16 = 2
    17 = 3
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 1
    24 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
